$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source/input variable for the "ID" row (row 2) changes from "D_ID" to "d_id"
$ws.Range("F2").Value = "d_id"

# Rows whose "input_variables" (column F) cell was empty and now gets the
# rule_category text "impossible " (trailing space) filled in, signalling
# that there is no source variable for these mappings.
$impossibleRows = @(18,19,22,24,29,32,36,37,39,47,48,49,56,59,60,61,62,63,64,70,78,79,80,85,86,88,89,91,92,99,100,104,107,108,109,110,111,112,113,114,115,116,117,118,119,120,121,122,123,124,125,126,127,128)

foreach ($r in $impossibleRows) {
    $ws.Range("F$r").Value = "impossible "
}

# Row 12 is the lone exception and gets "impossible" without the trailing space.
$ws.Range("F12").Value = "impossible"
